$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Cypher query text for cell B2 (replaces the "Programs/Arms/Cases..."
# summary-count query with a per-case detail query that also pulls in
# demographics). Built as a here-string so the embedded double quotes,
# backticks and newlines inside the Cypher text don't need escaping.
$newQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE ss.disease_subtype IN ["Paget's Disease"] and sf.grouped_recurrence_score IN ["11-15"]and d.tumor_size_group In ["(1,2]"] and tp.chemotherapy_regimen In ["Dose dense AC followed by a taxane"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
'@

$ws.Range("B2").Value = $newQuery

# The wrapped text is now taller, so the row grows to fit it.
$ws.Rows.Item(2).RowHeight = 348

# Selection moves from D4 to B4.
$ws.Range("B4").Select()
